# Rename TypeTest / TypeTestRef entities (and their sheets) to use their
# fully-qualified (package-prefixed) names, and update the package names
# themselves ("test" -> "org_molgenis_test", "molgenis" -> "org_molgenis").

$wb = $excel.ActiveWorkbook

$wsPackages   = $wb.Worksheets.Item("packages")
$wsEntities   = $wb.Worksheets.Item("entities")
$wsAttributes = $wb.Worksheets.Item("attributes")

# --- packages sheet -------------------------------------------------
# row3: package "molgenis" (parent "org") -> "org_molgenis"
$wsPackages.Range("A3").Value = "org_molgenis"
# row2 col C: parent package reference "molgenis" -> "org_molgenis"
$wsPackages.Range("C2").Value = "org_molgenis"
# row2: package "test" (parent "molgenis") -> "org_molgenis_test"
$wsPackages.Range("A2").Value = "org_molgenis_test"

# --- entities sheet ---------------------------------------------------
# TypeTest / TypeTestRef entities now live in package "org_molgenis_test"
$wsEntities.Range("B2").Value = "org_molgenis_test"
$wsEntities.Range("B3").Value = "org_molgenis_test"

# --- attributes sheet --------------------------------------------------
# entity column: rows referencing TypeTestRef / TypeTest
$wsAttributes.Range("B2:B3").Value = "org_molgenis_test_TypeTestRef"

# refEntity column: rows referencing TypeTestRef
$wsAttributes.Range("D10:D13").Value = "org_molgenis_test_TypeTestRef"
$wsAttributes.Range("D36:D37").Value = "org_molgenis_test_TypeTestRef"
$wsAttributes.Range("D42:D43").Value = "org_molgenis_test_TypeTestRef"
$wsAttributes.Range("D48").Value = "org_molgenis_test_TypeTestRef"

$wsAttributes.Range("B4:B50").Value = "org_molgenis_test_TypeTest"

# --- rename the data sheets themselves ---------------------------------
$wsTypeTest = $wb.Worksheets.Item("TypeTest")
$wsTypeTest.Name = "org_molgenis_test_TypeTest"

$wsTypeTestRef = $wb.Worksheets.Item("TypeTestRef")
$wsTypeTestRef.Name = "org_molgenis_test_TypeTestRef"

# --- column width tweaks ------------------------------------------------
# (ColumnWidth is stored internally with a +5/6 padding offset, so we back
# that off here to land on the desired stored width)
$wsPackages.Columns.Item(1).ColumnWidth = 16.666666666666668   # -> 17.5
$wsEntities.Columns.Item(2).ColumnWidth = 18.166666666666668   # -> 19
$wsAttributes.Columns.Item(2).ColumnWidth = 23.498697916666668 # -> ~24.33
$wsAttributes.Columns.Item(4).ColumnWidth = 22.166666666666668 # -> 23

# --- selection / active sheet tweaks ------------------------------------
$wsTypeTestRef.Select()
$wsTypeTestRef.Range("H41").Select()

$wsPackages.Select()
$wsPackages.Range("A2").Select()

$wsEntities.Select()
$wsEntities.Range("B3").Select()

$wsAttributes.Select()
$wsAttributes.Range("D48").Select()

$wsPackages.Select()
